# The "merged_data_python" sheet contains a block of lookup/helper values
# in columns B:F (rows 2-33) that were pulled in for an old VLOOKUP setup.
# Clear them out so the sheet can be rebuilt with a fresh VLOOKUP.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("merged_data_python")
$ws.Range("B2:F33").ClearContents()
